$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "35.304.30"; DNumeric = $False; E = "  +0.42%  " }
    @{ Row = 3; D = "1.885.31"; DNumeric = $False; E = "  -0.78%  " }
    @{ Row = 4; D = $null; DNumeric = $False; E = "  -0.71%  " }
    @{ Row = 5; D = "245.68"; DNumeric = $True; E = "  -2.87%  " }
    @{ Row = 6; D = $null; DNumeric = $False; E = "  -1.02%  " }
    @{ Row = 7; D = $null; DNumeric = $False; E = "  -0.76%  " }
    @{ Row = 8; D = "43.46"; DNumeric = $True; E = "  +6.28%  " }
    @{ Row = 9; D = $null; DNumeric = $False; E = "  -1.67%  " }
    @{ Row = 10; D = "53.58"; DNumeric = $True; E = "  +1.31%  " }
    @{ Row = 11; D = "0.0740"; DNumeric = $True; E = "  -1.51%  " }
    @{ Row = 12; D = "0.0972"; DNumeric = $True; E = "  -1.14%  " }
    @{ Row = 13; D = "13.27"; DNumeric = $True; E = "  +2.22%  " }
    @{ Row = 14; D = "2.159.40"; DNumeric = $False; E = "  -0.80%  " }
    @{ Row = 15; D = $null; DNumeric = $False; E = "  +2.66%  " }
    @{ Row = 16; D = $null; DNumeric = $False; E = "  -1.39%  " }
    @{ Row = 17; D = "1.886.95"; DNumeric = $False; E = "  -0.72%  " }
    @{ Row = 18; D = "35.427.89"; DNumeric = $False; E = $null }
    @{ Row = 19; D = "72.93"; DNumeric = $True; E = "  -0.97%  " }
    @{ Row = 20; D = $null; DNumeric = $False; E = "  -1.51%  " }
    @{ Row = 21; D = "244.31"; DNumeric = $True; E = "  +0.78%  " }
    @{ Row = 22; D = $null; DNumeric = $False; E = "  -1.54%  " }
    @{ Row = 23; D = $null; DNumeric = $False; E = "  -2.08%  " }
    @{ Row = 24; D = "2.65"; DNumeric = $True; E = "  +9.88%  " }
    @{ Row = 25; D = $null; DNumeric = $False; E = "  -0.71%  " }
    @{ Row = 26; D = $null; DNumeric = $False; E = "  -6.75%  " }
    @{ Row = 27; D = "166.01"; DNumeric = $True; E = "  -0.48%  " }
    @{ Row = 28; D = $null; DNumeric = $False; E = "  -0.94%  " }
    @{ Row = 29; D = "18.29"; DNumeric = $True; E = "  -1.02%  " }
    @{ Row = 30; D = $null; DNumeric = $False; E = "  -1.96%  " }
    @{ Row = 31; D = "4.128.44"; DNumeric = $False; E = "  +0.00%  " }
    @{ Row = 32; D = $null; DNumeric = $False; E = "  +10.61%  " }
    @{ Row = 33; D = "4.27"; DNumeric = $True; E = "  -1.02%  " }
    @{ Row = 34; D = "0.0583"; DNumeric = $True; E = "  -4.46%  " }
    @{ Row = 35; D = $null; DNumeric = $False; E = "  -1.40%  " }
    @{ Row = 36; D = $null; DNumeric = $False; E = "  -0.80%  " }
    @{ Row = 37; D = $null; DNumeric = $False; E = "  -11.69%  " }
    @{ Row = 38; D = "0.848"; DNumeric = $True; E = "  -0.48%  " }
    @{ Row = 39; D = $null; DNumeric = $False; E = "  -2.44%  " }
    @{ Row = 40; D = $null; DNumeric = $False; E = "  +7.12%  " }
    @{ Row = 41; D = $null; DNumeric = $False; E = "  +2.77%  " }
    @{ Row = 42; D = "17.22"; DNumeric = $True; E = "  +0.11%  " }
    @{ Row = 43; D = "96.53"; DNumeric = $True; E = "  -4.96%  " }
    @{ Row = 44; D = $null; DNumeric = $False; E = "  -2.17%  " }
    @{ Row = 45; D = "1.299.42"; DNumeric = $False; E = "  -1.34%  " }
    @{ Row = 46; D = $null; DNumeric = $False; E = "  -4.85%  " }
    @{ Row = 47; D = "0.0797"; DNumeric = $True; E = "  +7.89%  " }
    @{ Row = 48; D = "12.32"; DNumeric = $True; E = "  +4.05%  " }
    @{ Row = 49; D = $null; DNumeric = $False; E = "  -2.49%  " }
    @{ Row = 50; D = $null; DNumeric = $False; E = "  -0.83%  " }
    @{ Row = 51; D = $null; DNumeric = $False; E = "  -5.47%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$row")
        if ($u.DNumeric) {
            # Force text storage so numeric-looking strings (e.g. "245.68")
            # are not auto-converted to a number by Excel, then drop the
            # temporary text number-format again so the cell keeps its
            # original (unstyled) appearance.
            $dCell.NumberFormat = "@"
            $dCell.Value = $u.D
            $dCell.ClearFormats()
        } else {
            $dCell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
